$d = $word.ActiveDocument
$d.Content.Find.Execute("taught it the ESE program", $true, $false, $false, $false, $false, $true, 1, $false, "taught in the ESE program", 2)
